$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume values refreshed by the scraper run.
# Each cell is inline text (e.g. "292.74", "-0.99%"), so force a text
# number format before assigning, then clear the format again so the
# cell keeps its original (unstyled) look -- only the text content changes.
$updates = @{
    "D2" = "292.74"
    "D3" = "40.66"
    "E3" = "-0.99%"
    "D4" = "5.043"
    "E4" = "-3.04%"
    "E5" = "-3.39%"
    "D6" = "1.580"
    "E6" = "-3.20%"
    "D7" = "0.9249"
    "E7" = "1.11%"
    "E8" = "-0.41%"
    "D9" = "0.1169"
    "E9" = "-4.86%"
    "D10" = "0.1745"
    "E10" = "-4.14%"
    "D11" = "0.08718"
    "E11" = "-4.40%"
    "D12" = "0.04206"
    "E12" = "-1.34%"
    "D13" = "0.1054"
    "E13" = "0.33%"
    "D14" = "0.001276"
    "E14" = "1.41%"
    "D15" = "0.005981"
    "E15" = "2.21%"
    "D16" = "3.362"
    "E16" = "0.57%"
    "D17" = "4.323"
    "E17" = "0.64%"
    "D18" = "0.3348"
    "E18" = "0.37%"
    "D19" = "7.648"
    "E19" = "3.87%"
    "D20" = "0.1360"
    "E20" = "-1.58%"
    "D21" = "0.2822"
    "E21" = "4.04%"
    "D22" = "0.03869"
    "E22" = "-3.93%"
    "D23" = "0.001295"
    "E23" = "2.50%"
    "D24" = "0.003627"
    "E24" = "-17.12%"
    "E25" = "0.68%"
    "D26" = "0.0003750"
    "E26" = "-95.01%"
    "D38" = "0.02313"
    "E38" = "-7.36%"
    "D39" = "0.04998"
    "E39" = "-5.77%"
    "D40" = "0.007746"
    "E40" = "-1.20%"
    "D42" = "0.1275"
    "E42" = "-2.99%"
    "D43" = "0.007442"
    "E43" = "11.69%"
    "D44" = "0.007842"
    "E44" = "-2.21%"
    "D45" = "0.3178"
    "E45" = "3.67%"
    "D46" = "0.00006696"
    "E46" = "-0.67%"
    "E47" = "0.67%"
    "E48" = "-18.12%"
    "D49" = "0.004231"
    "E49" = "36.39%"
    "E50" = "0.67%"
    "D51" = "0.0002015"
    "E51" = "0.67%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
